$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.999.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "'1.860.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'311.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.5126"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.96%  "
$ws.Range("D8").Value = "'0.3809"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.08292"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.31%  "
$ws.Range("D10").Value = "'1.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").Value = "'6.202"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "'20.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "'1.859.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "'7.195"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "'1.003"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "'0.00001094"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "'90.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'0.06603"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "'17.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'6.013"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").Value = "'28.007.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("D24").Value = "'2.228"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").Value = "'2.565"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("D26").Value = "'2.074.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").Value = "'157.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'20.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "'124.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").Value = "'0.1062"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'1.034"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").Value = "'5.601"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "'3.596"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'9.588"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("D35").Value = "'0.06545"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "'0.02421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").Value = "'0.2171"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "'1.206"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").Value = "'0.6406"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").Value = "'11.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("D42").Value = "'4.880"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").Value = "'13.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "'1.273"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").Value = "'3.651"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "'1.976"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'1.204"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "'120.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "'79.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "'0.06827"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.86%  "
